$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parish rows to append (rows 3..6), matching the "id.xlsx" import sample.
# Columns: A=_id  B=IdxPar  C=IdxParCourt  D=IdxReg  E=NomParoisse  F=Stockage  G=Type
$ids      = @(1, 3, 4, 5)
$idxPar   = @(3000000000, 3030000000, 3080000000, 3090000000)
$idxParCt = @(3000, 3030, 3080, 3090)
$idxReg   = @(3000, 3000, 3000, 3000)
$noms     = @("Les Chamberonnes", "Le Mont-sur-Lausanne", "Ecublens – Saint-Sulpice", "Renens")
$stockage = @(1, 1, 1, 1)
$types    = @("R", "P", "P", "P")

# Populate column by column (matches the order new shared strings were first
# introduced: all of column E, then all of column G).
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 2).Value = $idxPar[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 3).Value = $idxParCt[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 4).Value = $idxReg[$i]
}
$ws.Range("E3").ClearFormats() | Out-Null
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 5).Value = $noms[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 6).Value = $stockage[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($i + 3, 7).Value = $types[$i]
}

# Match the "vertical top" alignment already used by the numeric columns on
# row 2 (columns A-D and F; columns E/G stay unstyled like E2/G2).
$ws.Range("A3:D6").VerticalAlignment = -4160
$ws.Range("F3:F6").VerticalAlignment = -4160

# Move the active selection, as recorded in the saved workbook.
$ws.Range("G4").Select() | Out-Null
